$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.183027744293213
$ws.Range("B1").Value = 2.550841093063354
$ws.Range("C1").Value = 9.410435676574707
$ws.Range("D1").Value = 2.09122633934021
$ws.Range("E1").Value = 1.219111084938049
